{"js": "// ------------------------------------------------------------------\n// 1) Insert a new \"Meta description\" paragraph right after the title\n//    (Heading1) paragraph at the top of the document.\n// ------------------------------------------------------------------\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nconst metaOoxml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover Dolphin Gold with Stellar Jackpots slot review and play for free. Enjoy 5 reels, 40 paylines, free spin features, and fixed jackpots.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nconst titleWholeRange = titlePara.getRange(Word.RangeLocation.whole);\ntitleWholeRange.insertOoxml(metaOoxml, Word.InsertLocation.after);\nawait context.sync();\n\n// ------------------------------------------------------------------\n// 2) At the end of the document: drop the bold \"title\" paragraph and\n//    swap the italic \"description\" paragraph's text for the new\n//    image-generation prompt (keeping its italic run/formatting).\n// ------------------------------------------------------------------\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs2.items.length;\nconst lastTitlePara = paragraphs2.items[count - 2];\nlastTitlePara.delete();\nawait context.sync();\n\nconst paragraphs3 = body.paragraphs;\nparagraphs3.load(\"items\");\nawait context.sync();\n\nconst lastPara = paragraphs3.items[paragraphs3.items.length - 1];\nconst lastRange = lastPara.getRange(Word.RangeLocation.whole);\nlastRange.insertText(\n  \"Prompt: Create a cartoon-style feature image for Dolphin Gold with Stellar Jackpots that features a happy Maya warrior wearing glasses. The image should have an underwater theme with the dolphin and gold elements incorporated into the background. The Maya warrior should be holding a treasure chest and smiling at the viewer. Use bright colors and bold lines to make the image pop and attract attention to the game's exciting features. The image should convey the idea of adventure and treasure while also showcasing the game's playful and enjoyable aspects.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# 1) Insert a new \"Meta description\" paragraph right after the title\n#    (Heading1) paragraph at the top of the document.\n# ------------------------------------------------------------------\n$titlePara = $d.Paragraphs(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs(2)\n$metaRng = $metaPara.Range\n$metaRng.Collapse(0)\n\n$metaXml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover Dolphin Gold with Stellar Jackpots slot review and play for free. Enjoy 5 reels, 40 paylines, free spin features, and fixed jackpots.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$metaRng.InsertXML($metaXml)\n\n# ------------------------------------------------------------------\n# 2) At the end of the document: drop the bold \"title\" paragraph and\n#    swap the italic \"description\" paragraph's text for the new\n#    image-generation prompt (keeping its italic run/formatting).\n# ------------------------------------------------------------------\n$count = $d.Paragraphs.Count\n$lastTitlePara = $d.Paragraphs($count - 1)\n$lastTitlePara.Range.Delete()\n\n$lastPara = $d.Paragraphs($d.Paragraphs.Count)\n$lastRng = $lastPara.Range\n$textOnly = $d.Range($lastRng.Start, $lastRng.End - 1)\n$textOnly.Text = \"Prompt: Create a cartoon-style feature image for Dolphin Gold with Stellar Jackpots that features a happy Maya warrior wearing glasses. The image should have an underwater theme with the dolphin and gold elements incorporated into the background. The Maya warrior should be holding a treasure chest and smiling at the viewer. Use bright colors and bold lines to make the image pop and attract attention to the game's exciting features. The image should convey the idea of adventure and treasure while also showcasing the game's playful and enjoyable aspects.\"\n"}
